$d = $word.ActiveDocument

$replacements = @(
    @("218×3=", "219×5="),
    @("763×2=", "214×8="),
    @("208×7=", "636×3="),
    @("406×7=", "464×3="),
    @("777×9=", "375×6="),
    @("863×8=", "219×2="),
    @("462×4=", "329×4="),
    @("438×2=", "865×2="),
    @("336×8=", "767×2="),
    @("311×9=", "811×2="),
    @("340×3=", "293×2="),
    @("719×7=", "846×5="),
    @("552×8=", "926×4="),
    @("822×7=", "362×6="),
    @("342×4=", "600×3="),
    @("871×5=", "975×8="),
    @("464×2=", "212×9="),
    @("432×9=", "960×8="),
    @("760×7=", "279×3="),
    @("550×9=", "786×8="),
    @("428×8=", "751×7="),
    @("630×5=", "381×7="),
    @("263×7=", "670×9="),
    @("320×3=", "582×9="),
    @("831×5=", "232×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
